$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 54, shifting existing rows 54-123 down to 55-124.
$ws.Rows.Item(54).Insert()

# Populate the newly inserted row 54 with the new weekly price record.
$ws.Cells.Item(54, 1).Value = 11
$ws.Cells.Item(54, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(54, 3).Value = "Bíobío"
$ws.Cells.Item(54, 4).Value = 45079
$ws.Cells.Item(54, 5).Value = 8
$ws.Cells.Item(54, 6).Value = 100112012
$ws.Cells.Item(54, 7).Value = "Espinaca"
$ws.Cells.Item(54, 8).Value = "Sin especificar"
$ws.Cells.Item(54, 9).Value = "Primera"
$ws.Cells.Item(54, 10).Value = 40
$ws.Cells.Item(54, 11).Value = 5000
$ws.Cells.Item(54, 12).Value = 5500
$ws.Cells.Item(54, 13).Value = 5250
$ws.Cells.Item(54, 14).Value = "$/cuna 10 kilos"
$ws.Cells.Item(54, 15).Value = "Región Metropolitana"
$ws.Cells.Item(54, 16).Value = 525
$ws.Cells.Item(54, 17).Value = 10
$ws.Cells.Item(54, 18).Value = "Hortaliza"
